$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notas")
$ws.Range("A2:O2").ClearContents()
